$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "19÷5=3, 4"
$t.Cell(1, 2).Range.Text = "91÷2=45, 1"
$t.Cell(1, 3).Range.Text = "32÷8=4, 0"
$t.Cell(1, 4).Range.Text = "77÷6=12, 5"
$t.Cell(1, 5).Range.Text = "43÷9=4, 7"
$t.Cell(5, 1).Range.Text = "34÷9=3, 7"
$t.Cell(5, 2).Range.Text = "19÷7=2, 5"
$t.Cell(5, 3).Range.Text = "56÷8=7, 0"
$t.Cell(5, 4).Range.Text = "59÷7=8, 3"
$t.Cell(5, 5).Range.Text = "92÷9=10, 2"
$t.Cell(9, 1).Range.Text = "92÷3=30, 2"
$t.Cell(9, 2).Range.Text = "82÷5=16, 2"
$t.Cell(9, 3).Range.Text = "79÷4=19, 3"
$t.Cell(9, 4).Range.Text = "37÷3=12, 1"
$t.Cell(9, 5).Range.Text = "45÷3=15, 0"
$t.Cell(13, 1).Range.Text = "71÷9=7, 8"
$t.Cell(13, 2).Range.Text = "55÷5=11, 0"
$t.Cell(13, 3).Range.Text = "15÷7=2, 1"
$t.Cell(13, 4).Range.Text = "95÷6=15, 5"
$t.Cell(13, 5).Range.Text = "12÷9=1, 3"
$t.Cell(17, 1).Range.Text = "62÷7=8, 6"
$t.Cell(17, 2).Range.Text = "82÷4=20, 2"
$t.Cell(17, 3).Range.Text = "61÷5=12, 1"
$t.Cell(17, 4).Range.Text = "75÷3=25, 0"
$t.Cell(17, 5).Range.Text = "98÷7=14, 0"
